# TODO Lists, GM Section Start, Reorg, Malboro
#
# This script reproduces (against the headless Excel COM-interop runtime) the
# changes described by the target OOXML diff:
#   - weapons: just move the cell selection/view state (no data change)
#   - armor_accessories: Winged Belt's Tribe changes from Yuke -> Selkie
#   - chalice_accessories: add a "Description" column (header + one row of
#     flavour text for Crystal Feather); becomes the active sheet
#   - edible_accessories: add a "Description" column, and two new rows
#     (Phoenix Down, Strange Liquid)
#   - monsters_rva and Sheet10 (and its backing table) are removed entirely
#   - workbook ends up with chalice_accessories as the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. weapons: no data change, just move the view selection to E10
# ---------------------------------------------------------------------------
$wsWeapons = $wb.Worksheets("weapons")
$wsWeapons.Range("E10").Select()

# ---------------------------------------------------------------------------
# 2. armor_accessories: Winged Belt is Selkie exclusive now, not Yuke
# ---------------------------------------------------------------------------
$wsArmorAcc = $wb.Worksheets("armor_accessories")
$wsArmorAcc.Range("C2").Value = "Selkie"
$wsArmorAcc.Range("C3").Select()

# ---------------------------------------------------------------------------
# 3. chalice_accessories: add Description column
# ---------------------------------------------------------------------------
$wsChalice = $wb.Worksheets("chalice_accessories")
$wsChalice.Range("C1").Value = "Description"
$wsChalice.Range("C2").Value = "Piece of bark of a Myrrh tree."
$wsChalice.Range("C3").Select()

# ---------------------------------------------------------------------------
# 4. edible_accessories: add Description column + two new items
# ---------------------------------------------------------------------------
$wsEdible = $wb.Worksheets("edible_accessories")
$wsEdible.Range("C1").Value = "Description"
$wsEdible.Range("C2").Value = "Leaf of a Myrrh tree, glowing with clean air."

$wsEdible.Range("A3").Value = "Phoenix Down"
$wsEdible.Range("B3").Value = "When you lose your last heart, this revives you at 1 heart."
$wsEdible.Range("C3").Value = "Magic flower. Looks like a flame and a feather."

$wsEdible.Range("A4").Value = "Strange Liquid"
$wsEdible.Range("B4").Value = "Heal all your hearts"
$wsEdible.Range("C4").Value = "Rare drink. Probably Myrrh dew and herbs."

$wsEdible.Range("C5").Select()

# ---------------------------------------------------------------------------
# 5. Remove the monsters_rva sheet and the leftover Sheet10 (with its table)
# ---------------------------------------------------------------------------
$wsSheet10 = $wb.Worksheets("Sheet10")
$wsSheet10.ListObjects("Import1").Delete()
$wb.Worksheets("Sheet10").Delete()
$wb.Worksheets("monsters_rva").Delete()

# ---------------------------------------------------------------------------
# 6. chalice_accessories ends up the active tab
# ---------------------------------------------------------------------------
$wsChalice.Activate()
